# Refresh the cryptocurrency price table (columns B-E, rows 2-51).
# Mirrors the upstream GitHub Actions data-refresh commit: prices / 1h
# volume percentages move, and a handful of rows swap rank position
# (name + link + price + volume all change together for those rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "31.334.23"
$ws.Range("E2").Value = "  +1.92%  "

# Row 3
$ws.Range("D3").Value = "1.996.91"
$ws.Range("E3").Value = "  +2.23%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.90%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "255.49"
$ws.Range("E5").Value = "  +1.30%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7413"
$ws.Range("E6").Value = "  +17.11%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  +0.77%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3398"
$ws.Range("E8").Value = "  +6.82%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "27.46"
$ws.Range("E9").Value = "  +8.83%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07139"
$ws.Range("E10").Value = "  +4.35%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8296"
$ws.Range("E11").Value = "  +2.06%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08152"
$ws.Range("E12").Value = "  +2.38%  "

# Row 13
$ws.Range("D13").Value = "2.000.20"
$ws.Range("E13").Value = "  +2.38%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "99.81"
$ws.Range("E14").Value = "  -1.31%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.592"
$ws.Range("E15").Value = "  +4.46%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.39"
$ws.Range("E16").Value = "  +11.39%  "

# Row 17
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "31.307.63"
$ws.Range("E17").Value = "  +1.76%  "

# Row 18
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "266.05"
$ws.Range("E18").Value = "  -4.30%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.100"
$ws.Range("E19").Value = "  +8.79%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000008253"
$ws.Range("E20").Value = "  +7.26%  "

# Row 21
$ws.Range("D21").Value = "2.261.21"
$ws.Range("E21").Value = "  +3.25%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.006"
$ws.Range("E22").Value = "  +0.84%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.008"
$ws.Range("E23").Value = "  +1.02%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.045"
$ws.Range("E24").Value = "  +6.36%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.02"
$ws.Range("E25").Value = "  +5.81%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.13"
$ws.Range("E26").Value = "  -1.13%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.82"
$ws.Range("E27").Value = "  +1.82%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.362"
$ws.Range("E28").Value = "  +11.44%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1328"
$ws.Range("E29").Value = "  +16.07%  "

# Row 30
$ws.Range("E30").Value = "  +3.76%  "

# Row 31
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.375"
$ws.Range("E31").Value = "  +1.56%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.629"
$ws.Range("E32").Value = "  +3.32%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.433"
$ws.Range("E33").Value = "  +2.20%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05356"
$ws.Range("E34").Value = "  +7.12%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.283"
$ws.Range("E35").Value = "  +7.86%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7845"
$ws.Range("E36").Value = "  +8.33%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.802"
$ws.Range("E37").Value = "  +3.05%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02003"
$ws.Range("E38").Value = "  +2.29%  "

# Row 39
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.903"
$ws.Range("E39").Value = "  -0.61%  "

# Row 40
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "84.34"
$ws.Range("E40").Value = "  +8.84%  "

# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.764"
$ws.Range("E41").Value = "  +4.36%  "

# Row 42
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4648"
$ws.Range("E42").Value = "  +1.81%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.125"
$ws.Range("E43").Value = "  +4.64%  "

# Row 44
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8596"
$ws.Range("E44").Value = "  +1.79%  "

# Row 45
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "105.07"
$ws.Range("E45").Value = "  +2.51%  "

# Row 46
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.005"
$ws.Range("E46").Value = "  +0.79%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.06"
$ws.Range("E47").Value = "  +0.26%  "

# Row 48
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.696"
$ws.Range("E48").Value = "  +5.60%  "

# Row 49
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.32"
$ws.Range("E49").Value = "  +3.80%  "

# Row 50
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.571"
$ws.Range("E50").Value = "  +10.43%  "

# Row 51
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4323"
$ws.Range("E51").Value = "  +3.67%  "
